# Generate Report for Handoff
# The status report moved the "69edf1bb-..." entry from "In Translation"
# down to "Ready for handoff" (new handoff timestamps), which re-sorts it
# below the "c2d03bec-..." / "48430eff-..." entries on every sheet.

$wb = $excel.ActiveWorkbook

function Set-CellAndLink($ws, $cellAddr, $value, $linkDisplay) {
    $ws.Range($cellAddr).Value = $value
    if ($linkDisplay -ne $null) {
        foreach ($hl in $ws.Hyperlinks) {
            $addr = $hl.Range.Address()
            if ($addr -eq ('$' + $cellAddr.Substring(0,1) + '$' + $cellAddr.Substring(1))) {
                $hl.TextToDisplay = $linkDisplay
            }
        }
    }
}

# ---------------------------------------------------------------------------
# Overview sheet: columns A (File Name, hyperlink), B (zh-cn), C (de-de),
# D (Latest Handoff Date)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

Set-CellAndLink $ws "A6" "c2d03bec-5599-426f-b692-644929f3c512.md" "c2d03bec-5599-426f-b692-644929f3c512.md"
$ws.Range("B6").Value = "In Translation"
$ws.Range("C6").Value = "In Translation"
$ws.Range("D6").Value = "2016-34-20 22:34:41"

Set-CellAndLink $ws "A7" "48430eff-6746-4dfc-b2d7-cbb467fa8e4c.md" "48430eff-6746-4dfc-b2d7-cbb467fa8e4c.md"
$ws.Range("B7").Value = "Ready for handoff"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("D7").Value = "2016-35-20 22:35:40"

Set-CellAndLink $ws "A8" "69edf1bb-0b0d-4ad6-aa21-988f17ab113b.md" "69edf1bb-0b0d-4ad6-aa21-988f17ab113b.md"
$ws.Range("B8").Value = "Ready for handoff"
$ws.Range("C8").Value = "Ready for handoff"
$ws.Range("D8").Value = "2016-38-20 22:38:39"

# ---------------------------------------------------------------------------
# zh-cn sheet: A (Source File Name), B (File Extension), C (Status),
# D (Latest Handoff File), E (Latest Handoff Datetime)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

Set-CellAndLink $ws "A6" "c2d03bec-5599-426f-b692-644929f3c512.md" "c2d03bec-5599-426f-b692-644929f3c512.md"
$ws.Range("C6").Value = "In Translation"
Set-CellAndLink $ws "D6" "c2d03bec-5599-426f-b692-644929f3c512.29f4087afcadea90589da1da82ef4653baabb6c9.zh-cn.xlf" "c2d03bec-5599-426f-b692-644929f3c512.29f4087afcadea90589da1da82ef4653baabb6c9.zh-cn.xlf"
$ws.Range("E6").Value = "2016-03-20 22:34:38"

Set-CellAndLink $ws "A7" "48430eff-6746-4dfc-b2d7-cbb467fa8e4c.md" "48430eff-6746-4dfc-b2d7-cbb467fa8e4c.md"
$ws.Range("C7").Value = "Ready for handoff"
Set-CellAndLink $ws "D7" "48430eff-6746-4dfc-b2d7-cbb467fa8e4c.b3673e249a8c96442fc9ec0006f263142f2c94fd.zh-cn.xlf" "48430eff-6746-4dfc-b2d7-cbb467fa8e4c.b3673e249a8c96442fc9ec0006f263142f2c94fd.zh-cn.xlf"
$ws.Range("E7").Value = "2016-03-20 22:35:37"

Set-CellAndLink $ws "A8" "69edf1bb-0b0d-4ad6-aa21-988f17ab113b.md" "69edf1bb-0b0d-4ad6-aa21-988f17ab113b.md"
$ws.Range("C8").Value = "Ready for handoff"
Set-CellAndLink $ws "D8" "69edf1bb-0b0d-4ad6-aa21-988f17ab113b.3816c468db1c35d1761145ddc234d2edf2c28616.zh-cn.xlf" "69edf1bb-0b0d-4ad6-aa21-988f17ab113b.3816c468db1c35d1761145ddc234d2edf2c28616.zh-cn.xlf"
$ws.Range("E8").Value = "2016-03-20 22:38:36"

# ---------------------------------------------------------------------------
# de-de sheet: same layout as zh-cn
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

Set-CellAndLink $ws "A6" "c2d03bec-5599-426f-b692-644929f3c512.md" "c2d03bec-5599-426f-b692-644929f3c512.md"
$ws.Range("C6").Value = "In Translation"
Set-CellAndLink $ws "D6" "c2d03bec-5599-426f-b692-644929f3c512.29f4087afcadea90589da1da82ef4653baabb6c9.de-de.xlf" "c2d03bec-5599-426f-b692-644929f3c512.29f4087afcadea90589da1da82ef4653baabb6c9.de-de.xlf"
$ws.Range("E6").Value = "2016-03-20 22:34:41"

Set-CellAndLink $ws "A7" "48430eff-6746-4dfc-b2d7-cbb467fa8e4c.md" "48430eff-6746-4dfc-b2d7-cbb467fa8e4c.md"
$ws.Range("C7").Value = "Ready for handoff"
Set-CellAndLink $ws "D7" "48430eff-6746-4dfc-b2d7-cbb467fa8e4c.b3673e249a8c96442fc9ec0006f263142f2c94fd.de-de.xlf" "48430eff-6746-4dfc-b2d7-cbb467fa8e4c.b3673e249a8c96442fc9ec0006f263142f2c94fd.de-de.xlf"
$ws.Range("E7").Value = "2016-03-20 22:35:40"

Set-CellAndLink $ws "A8" "69edf1bb-0b0d-4ad6-aa21-988f17ab113b.md" "69edf1bb-0b0d-4ad6-aa21-988f17ab113b.md"
$ws.Range("C8").Value = "Ready for handoff"
Set-CellAndLink $ws "D8" "69edf1bb-0b0d-4ad6-aa21-988f17ab113b.3816c468db1c35d1761145ddc234d2edf2c28616.de-de.xlf" "69edf1bb-0b0d-4ad6-aa21-988f17ab113b.3816c468db1c35d1761145ddc234d2edf2c28616.de-de.xlf"
$ws.Range("E8").Value = "2016-03-20 22:38:39"
